# The commit removes the speaker-notes text ("non-hierarchical" / "如何將Vit結合FPN")
# from slide 26's notes body placeholder, leaving the paragraph empty.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(26)
$notes = $s.NotesPage

for ($i = 1; $i -le $notes.Shapes.Count; $i++) {
    $shp = $notes.Shapes.Item($i)
    if ($shp.HasTextFrame) {
        $tr = $shp.TextFrame.TextRange
        if ($tr.Text.Contains("non-hierarchical")) {
            $tr.Text = ""
        }
    }
}
